$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the objective text in the shared-string-backed cells.
# Row order stays the same (A1/A2/A3) but the wording is shortened/renumbered:
#   A1: "example objective 1"  -> "obj1"
#   A2: "example objective 12" -> "obj2"
#   A3: "example objective 4"  -> "obj3"
$ws.Range("A1").Value = "obj1"
$ws.Range("A2").Value = "obj2"
$ws.Range("A3").Value = "obj3"

# Move the active selection to L9 (previously J23).
$ws.Range("L9").Select()
